$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("A1").Style

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '29.883.02'
$ws.Range('D2').Style = $defaultStyle
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value2 = '  +0.52%  '
$ws.Range('E2').Style = $defaultStyle

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '1.894.34'
$ws.Range('D3').Style = $defaultStyle
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value2 = '  +0.43%  '
$ws.Range('E3').Style = $defaultStyle

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '1.000'
$ws.Range('D4').Style = $defaultStyle
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value2 = '  -0.15%  '
$ws.Range('E4').Style = $defaultStyle

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '0.7826'
$ws.Range('D5').Style = $defaultStyle
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value2 = '  -0.01%  '
$ws.Range('E5').Style = $defaultStyle

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '243.81'
$ws.Range('D6').Style = $defaultStyle
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value2 = '  +1.02%  '
$ws.Range('E6').Style = $defaultStyle

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '1.000'
$ws.Range('D7').Style = $defaultStyle
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value2 = '  -0.05%  '
$ws.Range('E7').Style = $defaultStyle

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.3135'
$ws.Range('D8').Style = $defaultStyle
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value2 = '  -0.55%  '
$ws.Range('E8').Style = $defaultStyle

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value2 = '  +1.90%  '
$ws.Range('E9').Style = $defaultStyle

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.07291'
$ws.Range('D10').Style = $defaultStyle
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value2 = '  +4.55%  '
$ws.Range('E10').Style = $defaultStyle

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value2 = '  +0.68%  '
$ws.Range('E11').Style = $defaultStyle

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.7729'
$ws.Range('D12').Style = $defaultStyle
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value2 = '  +1.39%  '
$ws.Range('E12').Style = $defaultStyle

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '5.466'
$ws.Range('D13').Style = $defaultStyle
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value2 = '  +3.80%  '
$ws.Range('E13').Style = $defaultStyle

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '1.858.72'
$ws.Range('D14').Style = $defaultStyle
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value2 = '  -1.65%  '
$ws.Range('E14').Style = $defaultStyle

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '93.97'
$ws.Range('D15').Style = $defaultStyle
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value2 = '  +2.37%  '
$ws.Range('E15').Style = $defaultStyle

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '6.203'
$ws.Range('D16').Style = $defaultStyle
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value2 = '  +5.34%  '
$ws.Range('E16').Style = $defaultStyle

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '29.876.57'
$ws.Range('D17').Style = $defaultStyle
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value2 = '  +0.46%  '
$ws.Range('E17').Style = $defaultStyle

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '13.94'
$ws.Range('D18').Style = $defaultStyle
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value2 = '  +0.91%  '
$ws.Range('E18').Style = $defaultStyle

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '246.75'
$ws.Range('D19').Style = $defaultStyle
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value2 = '  +1.61%  '
$ws.Range('E19').Style = $defaultStyle

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '0.000007809'
$ws.Range('D20').Style = $defaultStyle
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value2 = '  +1.77%  '
$ws.Range('E20').Style = $defaultStyle

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value2 = '  -0.06%  '
$ws.Range('E21').Style = $defaultStyle

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '2.144.18'
$ws.Range('D22').Style = $defaultStyle
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value2 = '  +0.42%  '
$ws.Range('E22').Style = $defaultStyle

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '8.115'
$ws.Range('D23').Style = $defaultStyle
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value2 = '  +0.10%  '
$ws.Range('E23').Style = $defaultStyle

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value2 = '  -0.24%  '
$ws.Range('E24').Style = $defaultStyle

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '0.1596'
$ws.Range('D25').Style = $defaultStyle
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value2 = '  -3.41%  '
$ws.Range('E25').Style = $defaultStyle

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '9.430'
$ws.Range('D26').Style = $defaultStyle
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value2 = '  +1.87%  '
$ws.Range('E26').Style = $defaultStyle

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '164.28'
$ws.Range('D27').Style = $defaultStyle
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value2 = '  -0.38%  '
$ws.Range('E27').Style = $defaultStyle

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value2 = '  +0.80%  '
$ws.Range('E28').Style = $defaultStyle

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value2 = '2.018'
$ws.Range('D29').Style = $defaultStyle
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value2 = '  -0.96%  '
$ws.Range('E29').Style = $defaultStyle

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value2 = '  +3.08%  '
$ws.Range('E30').Style = $defaultStyle

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value2 = '  +0.85%  '
$ws.Range('E31').Style = $defaultStyle

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '4.474'
$ws.Range('D32').Style = $defaultStyle
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value2 = '  +2.43%  '
$ws.Range('E32').Style = $defaultStyle

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '0.05567'
$ws.Range('D33').Style = $defaultStyle
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value2 = '  -1.69%  '
$ws.Range('E33').Style = $defaultStyle

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '4.064'
$ws.Range('D34').Style = $defaultStyle
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value2 = '  +0.98%  '
$ws.Range('E34').Style = $defaultStyle

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '1.240'
$ws.Range('D35').Style = $defaultStyle
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value2 = '  -1.18%  '
$ws.Range('E35').Style = $defaultStyle

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '0.7532'
$ws.Range('D36').Style = $defaultStyle
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value2 = '  +2.97%  '
$ws.Range('E36').Style = $defaultStyle

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value2 = '  +0.08%  '
$ws.Range('E37').Style = $defaultStyle

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '2.686'
$ws.Range('D38').Style = $defaultStyle
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value2 = '  +1.77%  '
$ws.Range('E38').Style = $defaultStyle

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value2 = '0.01937'
$ws.Range('D39').Style = $defaultStyle
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value2 = '  +2.03%  '
$ws.Range('E39').Style = $defaultStyle

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '2.800'
$ws.Range('D40').Style = $defaultStyle
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value2 = '  +1.21%  '
$ws.Range('E40').Style = $defaultStyle

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '1.138.10'
$ws.Range('D41').Style = $defaultStyle
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value2 = '  +11.89%  '
$ws.Range('E41').Style = $defaultStyle

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '0.4463'
$ws.Range('D42').Style = $defaultStyle
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value2 = '  +1.97%  '
$ws.Range('E42').Style = $defaultStyle

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '74.13'
$ws.Range('D43').Style = $defaultStyle
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value2 = '  +2.86%  '
$ws.Range('E43').Style = $defaultStyle

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '5.963'
$ws.Range('D44').Style = $defaultStyle
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value2 = '  +2.90%  '
$ws.Range('E44').Style = $defaultStyle

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '0.8526'
$ws.Range('D45').Style = $defaultStyle
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value2 = '  +2.12%  '
$ws.Range('E45').Style = $defaultStyle

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '1.001'
$ws.Range('D46').Style = $defaultStyle
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value2 = '  +0.03%  '
$ws.Range('E46').Style = $defaultStyle

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '1.888'
$ws.Range('D47').Style = $defaultStyle
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value2 = '  +2.21%  '
$ws.Range('E47').Style = $defaultStyle

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '102.52'
$ws.Range('D48').Style = $defaultStyle
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value2 = '  +0.50%  '
$ws.Range('E48').Style = $defaultStyle

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '3.108'
$ws.Range('D49').Style = $defaultStyle
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value2 = '  +7.48%  '
$ws.Range('E49').Style = $defaultStyle

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '7.523'
$ws.Range('D50').Style = $defaultStyle
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value2 = '  +2.08%  '
$ws.Range('E50').Style = $defaultStyle

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '9.737'
$ws.Range('D51').Style = $defaultStyle
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value2 = '  -1.51%  '
$ws.Range('E51').Style = $defaultStyle
